$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the three new header cells (AD1:AF1), copying the bold/bordered
# header style from an existing header cell (A1) so the new headers share
# the same cellXf (s="1") instead of minting a new style.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill the season-record columns (Wins/Losses/Ties) for every data row.
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 30).Value = 83
    $ws.Cells.Item($r, 31).Value = 79
    $ws.Cells.Item($r, 32).Value = 1
}

Write-Output "done"
